$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": PORCELANATO sales for CHONTASI SIMBAÑA SILVIA JANETH (row 7) ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M7").Value = 846.1900000000001

# --- Sheet "VENTA MENSUAL": julio sales for CHONTASI SIMBAÑA SILVIA JANETH (row 7) and the TOTAL row (22) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F7").Value = 2445.56
$wsMensual.Range("F22").Value = 36598.88

# --- Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO row (16) and TOTAL row (19) ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D16").Value = 27385.83
$wsCumpl.Range("E16").Value = 16880.41
$wsCumpl.Range("F16").Value = 0.6186617611976983
$wsCumpl.Range("D19").Value = 36598.88
$wsCumpl.Range("E19").Value = 28779.11762291769
$wsCumpl.Range("F19").Value = 0.5598042358392847

# Column D on "CUMPLIMIENTO MENSUAL" grew from width 13 to 14 (Excel's stored OOXML width
# equals the ColumnWidth "characters" value plus ~0.8333 padding, so back it out here).
$wsCumpl.Columns.Item(4).ColumnWidth = 13.166666666666666
